$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from an existing header cell (H1) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for I2:J13
$data = @(
    @(3, 3),
    @(5, 6),
    @(10, 10),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(5, 6),
    @(8, 9),
    @(7, 8),
    @(10, 10),
    @(9, 9),
    @(9, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
